$wb = $excel.ActiveWorkbook

# --- Sheet: Formula Documentation (summary WBS overview) ---
$ws1 = $wb.Worksheets.Item("Formula Documentation")

$ws1.Range("A1").Value = 'INFORMATION TECHNOLOGY IMPLEMENTATION WBS - COMPREHENSIVE PROJECT PLAN'
$ws1.Range("A4").Value = 'This IT Implementation WBS contains 11 comprehensive phases with 178 detailed tasks,'
$ws1.Range("A10").Value = '3. IT Design and Planning - Requirements analysis, system design, documentation'
$ws1.Range("A11").Value = '4. IT Compliance and Regulatory - Compliance frameworks, regulatory approvals'
$ws1.Range("A12").Value = '5. IT Procurement and Contracting - Vendor selection, procurement, equipment'
$ws1.Range("A13").Value = '6. IT Implementation - Installation, configuration, integration, testing'
$ws1.Range("A14").Value = '7. IT Testing and Validation - System testing, performance validation, acceptance'
$ws1.Range("A15").Value = '8. IT Documentation and Training - Technical documentation, user training'
$ws1.Range("A16").Value = '9. IT Walkthrough and Handover - Deployment preparation, system handover'
$ws1.Range("A17").Value = '10. IT Project Close-out - Administrative closure, lessons learned, final reporting'
$ws1.Range("A18").Value = '11. IT Steady-State Support - Production support, maintenance, continuous improvement'
$ws1.Range("A20").Value = 'IT SPECIFIC FEATURES:'
$ws1.Range("A21").Value = '- Process analysis and reengineering
- Workflow optimization and automation
- Performance metrics and KPI implementation
- Change management and training
- Quality assurance and control
- Continuous improvement processes
- IT excellence initiatives'

# --- Sheet: Operational_Implementation_WBS (detailed WBS rows) ---
$ws2 = $wb.Worksheets.Item("Operational_Implementation_WBS")

$ws2.Range("C2").Value = 'IT Project Initiation [Deliverable: Project Initiation Report]'
$ws2.Range("C3").Value = '    IT Project Charter Development [Deliverable: Signed Project Charter]'
$ws2.Range("C5").Value = '        IT Requirements Definition [Deliverable: Requirements Document]'
$ws2.Range("K5").Value = 'IT Analyst'
$ws2.Range("C8").Value = '    IT Stakeholder Management [Deliverable: Stakeholder Management Plan]'
$ws2.Range("C13").Value = '    IT Initial Planning [Deliverable: Initial Project Plan]'
$ws2.Range("C18").Value = '    IT Project Initiation Approval [Deliverable: Signed Initiation Approval]'
$ws2.Range("C19").Value = 'IT Budget Planning and Management [Deliverable: Budget Management Report]'
$ws2.Range("C20").Value = '    IT Budget Development [Deliverable: Approved Budget Plan]'
$ws2.Range("C22").Value = '        IT Equipment Costs [Deliverable: Equipment Cost Report]'
$ws2.Range("K22").Value = 'IT Procurement Lead'
$ws2.Range("C25").Value = '    IT Budget Control and Tracking [Deliverable: Budget Control System]'
$ws2.Range("C29").Value = '    IT Budget Management Approval [Deliverable: Budget Management Approval]'
$ws2.Range("C30").Value = 'IT Design and Planning [Deliverable: Design Package]'
$ws2.Range("K30").Value = 'IT Solution Architect'
$ws2.Range("C31").Value = '    IT Requirements Analysis [Deliverable: Requirements Specification]'
$ws2.Range("C36").Value = '    IT System Design [Deliverable: System Design Document]'
$ws2.Range("C41").Value = '    IT Documentation [Deliverable: Design Documentation Package]'
$ws2.Range("C43").Value = '    IT Design Approval [Deliverable: Signed Design Approval]'
$ws2.Range("C44").Value = 'IT Compliance and Regulatory [Deliverable: Compliance Report]'
$ws2.Range("C45").Value = '    IT Compliance Assessment [Deliverable: Compliance Assessment Report]'
$ws2.Range("C51").Value = '    IT Regulatory Approvals [Deliverable: Regulatory Approval Documents]'
$ws2.Range("C57").Value = '    IT Compliance Documentation [Deliverable: Compliance Documentation]'
$ws2.Range("C59").Value = '    IT Compliance Approval [Deliverable: Signed Compliance Approval]'
$ws2.Range("C60").Value = 'IT Procurement and Contracting [Deliverable: Procurement Package]'
$ws2.Range("C61").Value = '    IT Vendor Selection [Deliverable: Selected Vendor List]'
$ws2.Range("C68").Value = '    IT Equipment Procurement [Deliverable: Procured Equipment]'
$ws2.Range("C73").Value = '    IT Procurement Documentation [Deliverable: Procurement Documentation]'
$ws2.Range("C74").Value = '    IT Procurement Completion [Deliverable: Procurement Completion Certificate]'
$ws2.Range("C75").Value = 'IT Implementation [Deliverable: Implemented System]'
$ws2.Range("C76").Value = '    IT Installation and Setup [Deliverable: Installed System]'
$ws2.Range("C83").Value = '    IT Integration and Configuration [Deliverable: Integrated System]'
$ws2.Range("C91").Value = '    IT Final Implementation [Deliverable: Final Implementation Report]'
$ws2.Range("C97").Value = '    IT Implementation Completion [Deliverable: Implementation Completion Certificate]'
$ws2.Range("C98").Value = 'IT Testing and Validation [Deliverable: Testing Report]'
$ws2.Range("C99").Value = '    IT System Testing [Deliverable: System Test Results]'
$ws2.Range("C106").Value = '    IT User Acceptance Testing [Deliverable: UAT Results]'
$ws2.Range("C112").Value = '    IT Final Validation [Deliverable: Final Validation Report]'
$ws2.Range("C115").Value = '    IT Testing Completion [Deliverable: Testing Completion Certificate]'
$ws2.Range("C116").Value = 'IT Documentation and Training [Deliverable: Documentation and Training Package]'
$ws2.Range("C117").Value = '    IT Technical Documentation [Deliverable: Technical Documentation Package]'
$ws2.Range("C123").Value = '    IT User Training [Deliverable: Trained Users]'
$ws2.Range("C129").Value = '    IT Documentation and Training Completion [Deliverable: Documentation and Training Completion Certificate]'
$ws2.Range("C130").Value = 'IT Walkthrough and Handover [Deliverable: Handover Package]'
$ws2.Range("C131").Value = '    IT Pre-Deployment Activities [Deliverable: Pre-Deployment Checklist]'
$ws2.Range("C137").Value = '    IT Go-Live Activities [Deliverable: Go-Live Report]'
$ws2.Range("C142").Value = '    IT System Handover [Deliverable: System Handover Package]'
$ws2.Range("C146").Value = '    IT Handover Completion [Deliverable: Handover Completion Certificate]'
$ws2.Range("C147").Value = 'IT Project Close-out [Deliverable: Project Close-out Report]'
$ws2.Range("C148").Value = '    IT Administrative Closure [Deliverable: Administrative Closure Package]'
$ws2.Range("C153").Value = '    IT Lessons Learned [Deliverable: Lessons Learned Report]'
$ws2.Range("C157").Value = '    IT Final Reporting [Deliverable: Final Project Report]'
$ws2.Range("C159").Value = '    IT Project Closure [Deliverable: Project Closure Certificate]'
$ws2.Range("C160").Value = 'IT Steady-State Support [Deliverable: Steady-State Support Plan]'
$ws2.Range("C161").Value = '    IT Production Support [Deliverable: Production Support Services]'
$ws2.Range("C167").Value = '    IT Maintenance and Updates [Deliverable: Maintenance Services]'
$ws2.Range("C172").Value = '    IT Continuous Improvement [Deliverable: Improvement Initiatives]'
$ws2.Range("C176").Value = '    IT Steady-State Establishment [Deliverable: Steady-State Establishment Certificate]'
